$d = $word.ActiveDocument

# 1) Mark every inline picture's range as "do not spell/grammar check"
#    (adds <w:noProof/> to the run properties of the runs hosting the
#    w:drawing elements for all 8 images in the document).
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = $true
}

# 2) Split the run that currently reads
#    " para actualizar las modificaciones del compañero anterior. Comando: "
#    into three runs:
#      " para actualizar las modificaciones del compañero anterior."
#      " "
#      "Comando: "
$old = " para actualizar las modificaciones del compañero anterior. Comando: "
$rng = $d.Content
$found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $part1 = " para actualizar las modificaciones del compañero anterior."
    $part2 = " "
    $part3 = "Comando: "

    $fullStart = $rng.Start

    $p1Start = $fullStart
    $p1End = $p1Start + $part1.Length

    $p2Start = $p1End
    $p2End = $p2Start + $part2.Length

    $p3Start = $p2End
    $p3End = $p3Start + $part3.Length

    $r2 = $d.Range($p2Start, $p2End)
    $r3 = $d.Range($p3Start, $p3End)

    # Toggling a character property forces the shared run to split into
    # separate runs at these boundaries; reverting the property back to
    # its default afterwards leaves the (now distinct) runs with
    # identical formatting to the original, but as separate <w:r>
    # elements in the saved XML.
    $r2.Font.Bold = $true
    $r3.Font.Bold = $true

    $r2b = $d.Range($p2Start, $p2End)
    $r3b = $d.Range($p3Start, $p3End)
    $r2b.Font.Bold = $false
    $r3b.Font.Bold = $false
}
